$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 "Marking" - Right/Wrong counts correction
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 "Total" - recompute totals to match corrected marking
$ws.Range("B12").Value = 64
$ws.Range("C12").Value = -4
$ws.Range("E12").Value = "60 / 112"
